$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All values in columns D (Price) and E (Volume label) are stored as text
# in the original workbook, so force the number format to Text ("@") before
# assigning, which prevents Excel from re-interpreting numeric-looking
# strings (e.g. "243.42") as actual numbers / floats.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "243.42"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "23.06"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.419"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05917"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.395"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8028"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9258"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1415"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07424"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.03395"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09347"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.941"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001595"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.04812"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0005944"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.005450"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.004330"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0009824"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.00007513"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.450"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1348"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03907"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006230"

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "40KickTokenKICK"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1074"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002615"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007285"

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "43LocalTradersLCTBestin24h"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005196"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0005804"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.002321"
